$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("C3").Value = "Circuitos Elétricos 2"
$ws.Range("D3").Value = "Circuitos Elétricos 2"
$ws.Range("F3").Value = "-"
$ws.Range("F6").Value = "EAP"
$ws.Range("F7").Value = "-"
